$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item 1s")

# B2 ("A Date" column, row for Id=0) now carries a time-of-day fraction
# instead of being a whole-day date serial.
$ws.Range("B2").Value = 41268.499652777777

# B3 ("A Date" column, row for Id=654) is truncated to a whole-day date
# serial (time-of-day fraction removed).
$ws.Range("B3").Value = 41275

# New row 4: Id=655 with a new String Field value ("Another string"),
# growing the used range to A1:C4.
$ws.Range("A4").Value = 655
$ws.Range("C4").Value = "Another string"
